$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-16 Thursday", "2025-01-17 Friday"),
    @("46×30=", "23×79="),
    @("42×60=", "15×83="),
    @("32×61=", "78×74="),
    @("88×29=", "57×13="),
    @("83×33=", "93×63="),
    @("24×62=", "80×24="),
    @("76×31=", "76×76="),
    @("45×43=", "86×92="),
    @("59×51=", "34×16="),
    @("35×25=", "78×33="),
    @("84×83=", "15×91="),
    @("50×47=", "96×16="),
    @("25×70=", "50×79="),
    @("67×45=", "80×45="),
    @("66×55=", "31×14="),
    @("12×81=", "15×15="),
    @("31×53=", "21×13="),
    @("68×55=", "41×33="),
    @("61×49=", "15×82="),
    @("47×81=", "32×89="),
    @("12×85=", "52×77="),
    @("74×76=", "35×40="),
    @("81×70=", "28×89="),
    @("13×22=", "15×88="),
    @("96×14=", "60×46=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
